$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  41"
$ws.Range("C9").Value = "Report Covering the Week  10/6/2025  Through  10/12/2025"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -95.588235294117
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 19
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = -5
$ws.Range("L15").Value = 35.714285714285
$ws.Range("M15").Value = 26.666666666666
$ws.Range("N15").Value = -68.852459016393
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -27.777777777777
$ws.Range("I16").Value = 144
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = -4
$ws.Range("L16").Value = -11.656441717791
$ws.Range("M16").Value = -24.607329842931
$ws.Range("N16").Value = -91.618160651920
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -27.272727272727
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -8.571428571428
$ws.Range("I17").Value = 366
$ws.Range("J17").Value = 304
$ws.Range("K17").Value = 20.394736842105
$ws.Range("L17").Value = 9.580838323353
$ws.Range("M17").Value = 48.178137651821
$ws.Range("N17").Value = -63.030303030303
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 117
$ws.Range("J18").Value = 118
$ws.Range("K18").Value = -0.847457627118
$ws.Range("L18").Value = -15.827338129496
$ws.Range("M18").Value = -51.851851851851
$ws.Range("N18").Value = -90.979182729375
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 18.75
$ws.Range("I19").Value = 338
$ws.Range("J19").Value = 285
$ws.Range("K19").Value = 18.596491228070
$ws.Range("L19").Value = 4.320987654320
$ws.Range("M19").Value = 45.064377682403
$ws.Range("N19").Value = -26.839826839826
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 93
$ws.Range("J20").Value = 89
$ws.Range("K20").Value = 4.494382022471
$ws.Range("L20").Value = -31.617647058823
$ws.Range("M20").Value = -8.823529411764
$ws.Range("N20").Value = -84.369747899159
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -25
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -12.711864406779
$ws.Range("I21").Value = 1080
$ws.Range("J21").Value = 976
$ws.Range("K21").Value = 10.655737704918
$ws.Range("L21").Value = -3.225806451612
$ws.Range("M21").Value = 3.250478011472
$ws.Range("N21").Value = -79.194760161818
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 20
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = 11.111111111111
$ws.Range("L22").Value = -9.090909090909
$ws.Range("M22").Value = -23.076923076923
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 11.111111111111
$ws.Range("I23").Value = 95
$ws.Range("J23").Value = 92
$ws.Range("K23").Value = 3.260869565217
$ws.Range("L23").Value = -2.061855670103
$ws.Range("M23").Value = 66.666666666666
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -20.833333333333
$ws.Range("F24").Value = 63
$ws.Range("H24").Value = -13.698630136986
$ws.Range("I24").Value = 669
$ws.Range("J24").Value = 570
$ws.Range("K24").Value = 17.368421052631
$ws.Range("L24").Value = 9.852216748768
$ws.Range("M24").Value = -0.149253731343
$ws.Range("C25").Value = 5
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 10
$ws.Range("I25").Value = 73
$ws.Range("K25").Value = -24.742268041237
$ws.Range("L25").Value = -8.75
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -57.142857142857
$ws.Range("F26").Value = 33
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = -28.260869565217
$ws.Range("I26").Value = 342
$ws.Range("J26").Value = 401
$ws.Range("K26").Value = -14.713216957606
$ws.Range("L26").Value = -23.146067415730
$ws.Range("M26").Value = -47.465437788018
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 26
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = -25.714285714285
$ws.Range("L27").Value = 8.333333333333
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 80
$ws.Range("I28").Value = 69
$ws.Range("K28").Value = 43.75
$ws.Range("L28").Value = 15
$ws.Range("I29").Value = 16
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = -20
$ws.Range("M29").Value = -74.193548387096
$ws.Range("N29").Value = -92.825112107623
$ws.Range("I30").Value = 13
$ws.Range("K30").Value = -40.909090909090
$ws.Range("L30").Value = -27.777777777777
$ws.Range("M30").Value = -72.916666666666
$ws.Range("N30").Value = -93.532338308457
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 3
$ws.Range("K31").Value = 50
$ws.Range("L31").Value = 0

# --- Cells changing from text placeholder to numeric (restyle to numeric format) ---
$ws.Range("C15").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = 0
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = 4
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = 2
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("C29").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("C30").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122) | Out-Null
$ws.Range("F31").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null

# --- Cells changing from numeric to text placeholder (restyle to general/text format) ---
$ws.Range("D25").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
